$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells for Wins / Losses / Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
